# CardList.xlsx update: rename Sheet1 -> Master, add Tier 1/2/3 sheets,
# add a design-note callout box on Master, and build out ability tables
# on the new Tier sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Rename the existing sheet to "Master"
# ---------------------------------------------------------------
$master = $wb.Worksheets.Item(1)
$master.Name = "Master"

# ---------------------------------------------------------------
# 2) Add the three new "Tier" sheets after Master
# ---------------------------------------------------------------
$tier1 = $wb.Worksheets.Add($null, $master)
$tier1.Name = "Tier 1"

$tier2 = $wb.Worksheets.Add($null, $tier1)
$tier2.Name = "Tier 2"

$tier3 = $wb.Worksheets.Add($null, $tier2)
$tier3.Name = "Tier 3"

Write-Host "Sheets created."
